$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers ("573.51", "8.00", ...).
# The source workbook stores Price-column values as TEXT (not Number) so
# exact formatting (e.g. trailing zeros such as "8.00") is preserved.
# Force NumberFormat to Text ("@") before assigning so Excel does not
# auto-convert the literal into a numeric value.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

$ws.Range("D2").Value = "62.796.12"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.461.28"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "573.51"
$ws.Range("D6").Value = "146.22"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "2.461.22"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "28.99"
$ws.Range("E14").Value = "  +2.73%  "
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").Value = "2.907.13"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "62.682.84"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "2.457.65"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "8.00"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "326.65"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("E22").Value = "  +9.89%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "10.20"
$ws.Range("E25").Value = "  +20.72%  "
$ws.Range("D26").Value = "65.58"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").Value = "659.26"
$ws.Range("E27").Value = "  +3.04%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0977"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.580.60"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -14.85%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("E32").Value = "  -2.47%  "
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("D37").Value = "4.75"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "0.369"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("D40").Value = "151.47"
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").Value = "18.70"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("D44").Value = "0.0₆0311"
$ws.Range("E44").Value = "  -65.74%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "152.83"
$ws.Range("E46").Value = "  +5.21%  "
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("E51").Value = "  -0.94%  "
